$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 4 new rows into the data table (before the former last row, 25)
#    so the "closing" styled row (old row 25) shifts down to row 29, and the
#    footer block (old rows 30:31) shifts down to rows 34:35, all via a
#    native Excel row-insert so formatting / merges shift automatically.
# ---------------------------------------------------------------------------
$ws.Rows("25:28").Insert()

# Re-apply the standard data-row formatting (copied from row 24, the last of
# the original "regular" styled rows) onto the 4 freshly inserted rows.
$ws.Range("B24:J24").Copy() | Out-Null
$ws.Range("B25:J28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Update the small header / summary cells
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "ESTADO DE CUENTA"
$ws.Range("B7").Value = "RAZON SOCIAL:"
$ws.Range("B11").Value = "VALOR MORA"
$ws.Range("E11").Value = 779316
$ws.Range("B13").Value = "Cant. Trabajadores"
$ws.Range("E13").Value = "Cant. Periodos"
$ws.Range("F13").Value = 6
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"
$ws.Range("J15").Value = "Observaciones"

# ---------------------------------------------------------------------------
# 3) Rewrite the 14-row data table (rows 16-29)
# ---------------------------------------------------------------------------
$data = @(
    @("CC","73165996","WILSON ENRIQUE LORDUY LLERENA","2507",3796,1300000),
    @("CC","73165996","WILSON ENRIQUE LORDUY LLERENA","2505",56940,1300000),
    @("CC","73165996","WILSON ENRIQUE LORDUY LLERENA","2503",56940,1300000),
    @("CC","73165996","WILSON ENRIQUE LORDUY LLERENA","2502",56940,1300000),
    @("CC","1143352761","ANDREA PAOLA VILLALOBOS SIMANCAS","2507",64000,1600000),
    @("CC","1143352761","ANDREA PAOLA VILLALOBOS SIMANCAS","2506",64000,1600000),
    @("CC","1143352761","ANDREA PAOLA VILLALOBOS SIMANCAS","2505",64000,1600000),
    @("CC","1143352761","ANDREA PAOLA VILLALOBOS SIMANCAS","2504",64000,1600000),
    @("CC","1143352761","ANDREA PAOLA VILLALOBOS SIMANCAS","2503",64000,1600000),
    @("CC","9294312","JHON JAIRO PAJARO ROJANO","2507",56940,1423500),
    @("CC","9294312","JHON JAIRO PAJARO ROJANO","2506",56940,1423500),
    @("CC","9294312","JHON JAIRO PAJARO ROJANO","2505",56940,1423500),
    @("CC","9294312","JHON JAIRO PAJARO ROJANO","2504",56940,1423500),
    @("CC","9294312","JHON JAIRO PAJARO ROJANO","2503",56940,1423500)
)

$r = 16
foreach ($row in $data) {
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $r++
}

Write-Host "Done"
